# Auto update stock data
# Updates the "as of" row (row 2) on each ticker sheet with refreshed
# date / EBITDA / ratio figures, and refreshes the Altman Z-Score /
# Piotroski F-Score columns where the underlying scores changed.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $text) {
    # Force the cell to stay text (matches the workbook's inlineStr cells)
    # instead of letting Excel auto-convert look-alike dates/numbers.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Alcoa (AA) ---
$ws = $wb.Worksheets.Item("Alcoa")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "4.79"

# --- Rio Tinto (RIO) ---
$ws = $wb.Worksheets.Item("Rio Tinto")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "7.59"

# --- Norsk Hydro (NHY) ---
$ws = $wb.Worksheets.Item("Norsk Hydro")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "2.69"
$ws.Range("G2:G8").Value = 3.25

# --- Reliance Steel & Aluminum (RS) ---
$ws = $wb.Worksheets.Item("Reliance Steel & Aluminum")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "12.27"

# --- Kaiser Aluminum (KALU) ---
$ws = $wb.Worksheets.Item("Kaiser Aluminum")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "9.96"
Set-TextValue $ws "C2" "1.33"
$ws.Range("H2:H8").Value = 6

# --- Ryerson Holding (RYI) ---
$ws = $wb.Worksheets.Item("Ryerson Holding")
Set-TextValue $ws "A2" "2025/10/25"
Set-TextValue $ws "B2" "20.57"
